$wb = $excel.ActiveWorkbook

# --- Sheet "Horario" (schedule grid) ---
$horario = $wb.Worksheets.Item("Horario")

# New course block: Algoritmos y Complejidad (3) now meets Miércoles 10:00-13:00 @ Ciencias 506
$horario.Range("D3").Value = "Algoritmos y Complejidad (3)`nCiencias 506"

# Introducción a la Programación (asdf) now meets Lunes 11:00-13:00 @ Reloj 102
$horario.Range("B4").Value = "Introducción a la Programación (asdf)`nReloj 102"

# Algoritmos y Complejidad (3) + new Diseño de Software Verificable (A), Miércoles 11:00-13:00
$horario.Range("D4").Value = "Algoritmos y Complejidad (3)`nCiencias 506`nDiseño de Software Verificable (A)`nHumanidades 203"

$horario.Range("B5").Value = "Introducción a la Programación (asdf)`nReloj 102"
$horario.Range("D5").Value = "Algoritmos y Complejidad (3)`nCiencias 506`nDiseño de Software Verificable (A)`nHumanidades 203"

# Clear the old 14:00-16:00 block (Algoritmos moved, asdf section moved)
$horario.Range("B7").Value = ""
$horario.Range("B8").Value = ""
$horario.Range("F8").Value = ""
$horario.Range("B9").Value = ""
$horario.Range("F9").Value = ""

# Keep row heights at their original auto size (avoid introducing explicit
# custom row heights as a side effect of wrapping the new multi-line text)
$horario.Rows.Item(3).AutoFit()
$horario.Rows.Item(4).AutoFit()
$horario.Rows.Item(5).AutoFit()
$horario.Rows.Item(7).AutoFit()
$horario.Rows.Item(8).AutoFit()
$horario.Rows.Item(9).AutoFit()

# --- Sheet "Tabla" (flat course table) ---
$tabla = $wb.Worksheets.Item("Tabla")

# Update "Algoritmos y Complejidad (Sección 3)" row: Miércoles 10:00-13:00 @ Ciencias 506
$tabla.Range("B4").Value = "Miércoles"
$tabla.Range("C4").Value = "10:00"
$tabla.Range("D4").Value = "13:00"
$tabla.Range("E4").Value = "Ciencias 506"

# Update "Introducción a la Programación (Sección asdf)" row: Lunes 11:00-13:00 @ Reloj 102
$tabla.Range("B5").Value = "Lunes"
$tabla.Range("C5").Value = "11:00"
$tabla.Range("D5").Value = "13:00"
$tabla.Range("E5").Value = "Reloj 102"

# New row: "Diseño de Software Verificable (Sección A)", Miércoles 11:00-13:00 @ Humanidades 203
$tabla.Range("A6").Value = "Diseño de Software Verificable (Sección A)"
$tabla.Range("B6").Value = "Miércoles"
$tabla.Range("C6").Value = "11:00"
$tabla.Range("D6").Value = "13:00"
$tabla.Range("E6").Value = "Humanidades 203"
